# Weekly data refresh for Fruta/Hortaliza sheet (Chirimoya - Terminal Hortofruticola Agro Chillan).
# The underlying source table got re-pulled; most rows 2-14 and 17 shifted to a new
# week's figures (date, quality, volume, min/max/avg price, price per kg).
# Rows 15 and 16 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44446
$ws.Range("L2").Value = "Primera"
$ws.Range("N2").Value = 21000
$ws.Range("O2").Value = 22000
$ws.Range("P2").Value = 21500
$ws.Range("S2").Value = 2150

# Row 3
$ws.Range("D3").Value = 44447
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("S3").Value = 2150

# Row 4
$ws.Range("D4").Value = 44448
$ws.Range("L4").Value = "Primera"
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 2150

# Row 5
$ws.Range("D5").Value = 44460
$ws.Range("L5").Value = "Especial"
$ws.Range("N5").Value = 31000
$ws.Range("O5").Value = 32000
$ws.Range("P5").Value = 31500
$ws.Range("S5").Value = 3150

# Row 6
$ws.Range("D6").Value = 44460
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("S6").Value = 3000

# Row 7
$ws.Range("D7").Value = 44461
$ws.Range("L7").Value = "Especial"
$ws.Range("N7").Value = 31000
$ws.Range("O7").Value = 32000
$ws.Range("P7").Value = 31500
$ws.Range("S7").Value = 3150

# Row 8
$ws.Range("D8").Value = 44461
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 30000
$ws.Range("O8").Value = 30000
$ws.Range("P8").Value = 30000
$ws.Range("S8").Value = 3000

# Row 9
$ws.Range("D9").Value = 44452
$ws.Range("L9").Value = "Primera"
$ws.Range("N9").Value = 21000
$ws.Range("O9").Value = 22000
$ws.Range("P9").Value = 21500
$ws.Range("S9").Value = 2150

# Row 10
$ws.Range("D10").Value = 44841
$ws.Range("M10").Value = 60
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 23500
$ws.Range("S10").Value = 2350

# Row 11
$ws.Range("D11").Value = 44839
$ws.Range("M11").Value = 120
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 26000
$ws.Range("P11").Value = 25500
$ws.Range("S11").Value = 2550

# Row 12
$ws.Range("D12").Value = 44487
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 23000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 23500
$ws.Range("S12").Value = 2350

# Row 13
$ws.Range("D13").Value = 44848
$ws.Range("L13").Value = "Especial"
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 24500
$ws.Range("S13").Value = 2450

# Row 14
$ws.Range("D14").Value = 44848
$ws.Range("M14").Value = 120
$ws.Range("N14").Value = 21000
$ws.Range("O14").Value = 22000
$ws.Range("P14").Value = 21500
$ws.Range("S14").Value = 2150

# Row 17
$ws.Range("D17").Value = 44868
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 26000
$ws.Range("P17").Value = 26000
$ws.Range("S17").Value = 2600
